$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.114.31"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.780.92"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "225.61"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D8").Value = "31.99"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "2.038.45"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "1.783.53"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "10.89"
$ws.Range("E14").Value = "  -5.58%  "
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "34.084.27"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "67.52"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "245.30"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "10.86"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "162.09"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "16.25"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").Value = "1.444.57"
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("D36").Value = "2.48"
$ws.Range("E36").Value = "  +5.65%  "
$ws.Range("D37").Value = "0.649"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.0190"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("D40").Value = "81.15"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "13.63"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "0.0519"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "1.938.80"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -6.16%  "
$ws.Range("D50").Value = "104.68"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("E51").Value = "  +0.32%  "
